$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update column C (Förändrad) for rows 2 through 12 from 45204 to 45207
$ws.Range("C2:C12").Value = 45207
